# Removed references to removed fields
#
# The "crop_name" calculate field and its accompanying "crop_disp" note
# (rows 16-18 of the "survey" sheet) referenced a removed ${crop_most}
# field and are no longer needed, so delete those rows outright. All
# subsequent rows (and their shared-string references) shift up to fill
# the gap, which Excel/the workbook engine handles automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Range("A16:A18").EntireRow.Delete()
